# "scroll stopped in experiment block, showseession fixes"
#
# Two things changed in the authored commit:
#  1) The scroll/selection position on the "interruptionsExperiment" sheet
#     (the sheet that was active/showing when the session was captured)
#     had stopped at N21; it should be parked at H23 instead.
#  2) The volatile experiment-randomisation block on "Sheet1" (A1:G9,
#     using RAND()/FLOOR()/IF(RAND()>0.5,...)) got recalculated, so its
#     cached values moved to a fresh set of random draws.

$wb = $excel.ActiveWorkbook

# 1) Fix the stuck scroll/selection on the interruptionsExperiment sheet.
$wsExp = $wb.Worksheets.Item("interruptionsExperiment")
$wsExp.Activate()
$wsExp.Range("H23").Select()
